$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Recipe-quantities table (row 1): "7500g" -> "750g" and
#    "1000g" -> "500g" (2nd & 3rd columns of the first data row).
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

# Column 2: "75" + "00g" ("7500g") -> "75" + "0g" ("750g")
$cell2 = $tbl.Cell(1, 2)
$c2Start = $cell2.Range.Start
$run2 = $d.Range($c2Start + 2, $c2Start + 5)
$r2 = $run2.Find.Execute("00g", $true, $false, $false, $false, $false,
                          $true, 0, $false, "0g", 1)

# Column 3: "1000g" -> "500g"
$cell3 = $tbl.Cell(1, 3)
$c3Start = $cell3.Range.Start
$run3 = $d.Range($c3Start, $c3Start + 5)
$r3 = $run3.Find.Execute("1000g", $true, $false, $false, $false, $false,
                          $true, 0, $false, "500g", 1)

# ------------------------------------------------------------------
# 2) Add a new instruction paragraph right before
#    "Ajouter les ingrédients dans l'ordre."
# ------------------------------------------------------------------
$paraCount = $d.Content.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Content.Paragraphs.Item($i)
    if ($p.Range.Text -match "Ajouter les ingrédients dans l.ordre") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -ge 1) {
    $target = $d.Content.Paragraphs.Item($targetIdx)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Content.Paragraphs.Item($targetIdx)
    $newPara.Range.Text = "Clipser le mélangeur de la machine."
}

# ------------------------------------------------------------------
# 3) Typo fix: "La cuisson prends environ " -> "La cuisson prend environ "
# ------------------------------------------------------------------
$content = $d.Content
$found = $content.Find.Execute("La cuisson prends environ ", $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
if ($found) {
    $fixStart = $content.Start
    $fixEnd = $content.End
    $fixRng = $d.Range($fixStart, $fixEnd)
    $r4 = $fixRng.Find.Execute("La cuisson prends environ ", $true, $false, $false, $false, $false,
                                $true, 0, $false, "La cuisson prend environ ", 1)
}

Write-Host "Edits applied: qty2=$r2 qty3=$r3 newPara=$($targetIdx -ge 1) typo=$r4"
